$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H19").Value = 599.62964
$ws.Range("I19").Value = 495.85715
$ws.Range("J19").Value = 635.95
$ws.Range("K19").Value = 495.85715
$ws.Range("L19").Value = 635.95
$ws.Range("M19").Value = -320.85715
$ws.Range("N19").Value = -985.95
$ws.Range("H86").Value = 1229.8889
$ws.Range("I86").Value = 1229.8889
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1229.8889
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -106.8888999999999
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1229.8889
$ws.Range("I89").Value = 1229.8889
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6149.4445
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -533.4444999999996
$ws.Range("N89").ClearContents()
$ws.Range("H129").Value = 994.65625
$ws.Range("J129").Value = 1096.8928
$ws.Range("L129").Value = 3290.6784
$ws.Range("N129").Value = -13290.6784
$ws.Range("H137").Value = 1638.28
$ws.Range("I137").Value = 1023.5
$ws.Range("J137").Value = 1755.381
$ws.Range("K137").Value = 3070.5
$ws.Range("L137").Value = 5266.143
$ws.Range("M137").Value = -520.5
$ws.Range("N137").Value = -10366.143
$ws.Range("H138").Value = 31254184
$ws.Range("I138").Value = 3912141
$ws.Range("J138").Value = 55558224
$ws.Range("K138").Value = 11736423
$ws.Range("L138").Value = 166674672
$ws.Range("M138").Value = -11731283
$ws.Range("N138").Value = -166684952

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H32").Value = 1811.1578
$ws.Range("I32").Value = 1281.5842
$ws.Range("K32").Value = 1281.5842
$ws.Range("M32").Value = -994.5842
$ws.Range("H61").Value = 1680.3928
$ws.Range("I61").Value = 1005
$ws.Range("J61").Value = 5732.75
$ws.Range("K61").Value = 1005
$ws.Range("L61").Value = 5732.75
$ws.Range("M61").Value = -793
$ws.Range("N61").Value = -6156.75
$ws.Range("H74").Value = 5364.839
$ws.Range("I74").Value = 1626.6666
$ws.Range("J74").Value = 13215
$ws.Range("K74").Value = 1626.6666
$ws.Range("L74").Value = 13215
$ws.Range("M74").Value = -752.6666
$ws.Range("N74").Value = -14963
$ws.Range("H77").Value = 5364.839
$ws.Range("I77").Value = 1626.6666
$ws.Range("J77").Value = 13215
$ws.Range("K77").Value = 8133.333000000001
$ws.Range("L77").Value = 66075
$ws.Range("M77").Value = -3765.333000000001
$ws.Range("N77").Value = -74811
$ws.Range("H122").Value = 2326
$ws.Range("I122").Value = 1812
$ws.Range("J122").Value = 2428.8
$ws.Range("K122").Value = 5436
$ws.Range("L122").Value = 7286.400000000001
$ws.Range("M122").Value = -2986
$ws.Range("N122").Value = -12186.4
$ws.Range("H132").Value = 1769.9803
$ws.Range("I132").Value = 1422.3695
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 4267.1085
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -1737.1085
$ws.Range("N132").Value = -19964
$ws.Range("H136").Value = 1680.3928
$ws.Range("I136").Value = 1005
$ws.Range("J136").Value = 5732.75
$ws.Range("K136").Value = 3015
$ws.Range("L136").Value = 17198.25
$ws.Range("M136").Value = -465
$ws.Range("N136").Value = -22298.25
$ws.Range("H139").Value = 60190.445
$ws.Range("J139").Value = 60190.445
$ws.Range("L139").Value = 60190.445
$ws.Range("N139").Value = -70470.44500000001

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 2340.484
$ws.Range("I134").Value = 1330.1666
$ws.Range("J134").Value = 5804.4287
$ws.Range("K134").Value = 3990.4998
$ws.Range("L134").Value = 17413.2861
$ws.Range("M134").Value = -1455.4998
$ws.Range("N134").Value = -22483.2861

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 161
$ws.Range("I16").Value = 161
$ws.Range("K16").Value = 161
$ws.Range("M16").Value = 126
$ws.Range("H62").Value = 21655
$ws.Range("I62").Value = 30457.857
$ws.Range("J62").Value = 6250
$ws.Range("K62").Value = 30457.857
$ws.Range("L62").Value = 6250
$ws.Range("M62").Value = -29833.857
$ws.Range("N62").Value = -7498
$ws.Range("H65").Value = 21655
$ws.Range("I65").Value = 30457.857
$ws.Range("J65").Value = 6250
$ws.Range("K65").Value = 152289.285
$ws.Range("L65").Value = 31250
$ws.Range("M65").Value = -149169.285
$ws.Range("N65").Value = -37490
$ws.Range("H113").Value = 161
$ws.Range("I113").Value = 161
$ws.Range("K113").Value = 161
$ws.Range("M113").Value = 2009

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H2").Value = 52631624
$ws.Range("I2").Value = 33.76923
$ws.Range("K2").Value = 202.61538
$ws.Range("M2").Value = -89.61538000000002
$ws.Range("H68").Value = 1341.49
$ws.Range("I68").Value = 974.7222
$ws.Range("J68").Value = 1772.0435
$ws.Range("K68").Value = 2924.1666
$ws.Range("L68").Value = 5316.1305
$ws.Range("M68").Value = -2113.1666
$ws.Range("N68").Value = -6938.1305
$ws.Range("H71").Value = 1341.49
$ws.Range("I71").Value = 974.7222
$ws.Range("J71").Value = 1772.0435
$ws.Range("K71").Value = 8772.4998
$ws.Range("L71").Value = 15948.3915
$ws.Range("M71").Value = -4716.4998
$ws.Range("N71").Value = -24060.3915
$ws.Range("H107").Value = 1618.2029
$ws.Range("I107").Value = 663.75
$ws.Range("J107").Value = 1743.3771
$ws.Range("K107").Value = 1991.25
$ws.Range("L107").Value = 5230.1313
$ws.Range("M107").Value = -71.25
$ws.Range("N107").Value = -9070.131300000001
$ws.Range("H113").Value = 11905201
$ws.Range("I113").Value = 365.33334
$ws.Range("J113").Value = 12820958
$ws.Range("K113").Value = 1096.00002
$ws.Range("L113").Value = 38462874
$ws.Range("M113").Value = 1073.99998
$ws.Range("N113").Value = -38467214
$ws.Range("H131").Value = 3048.7795
$ws.Range("J131").Value = 3341.698
$ws.Range("L131").Value = 10025.094
$ws.Range("N131").Value = -20105.094

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H122").Value = 2780152.8
$ws.Range("I122").Value = 3706203.8
$ws.Range("K122").Value = 11118611.4
$ws.Range("M122").Value = -11116161.4
$ws.Range("H132").Value = 2610.9167
$ws.Range("I132").Value = 2237.12
$ws.Range("J132").Value = 3460.4546
$ws.Range("K132").Value = 6711.36
$ws.Range("L132").Value = 10381.3638
$ws.Range("M132").Value = -4181.36
$ws.Range("N132").Value = -15441.3638
$ws.Range("H133").Value = 12400
$ws.Range("J133").Value = 12400
$ws.Range("L133").Value = 12400
$ws.Range("N133").Value = -22520
$ws.Range("H138").Value = 65500
$ws.Range("J138").Value = 65500
$ws.Range("L138").Value = 65500
$ws.Range("N138").Value = -75780

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H16").Value = 1030.3334
$ws.Range("I16").Value = 945.5
$ws.Range("K16").Value = 945.5
$ws.Range("M16").Value = -775.5
$ws.Range("H61").Value = 4374.316
$ws.Range("I61").Value = 5361.6924
$ws.Range("J61").Value = 2235
$ws.Range("K61").Value = 5361.6924
$ws.Range("L61").Value = 2235
$ws.Range("M61").Value = -5159.6924
$ws.Range("N61").Value = -2639
$ws.Range("H113").Value = 4374.316
$ws.Range("I113").Value = 5361.6924
$ws.Range("J113").Value = 2235
$ws.Range("K113").Value = 5361.6924
$ws.Range("L113").Value = 2235
$ws.Range("M113").Value = -3191.6924
$ws.Range("N113").Value = -6575
$ws.Range("H132").Value = 3124.4119
$ws.Range("I132").Value = 1412.5714
$ws.Range("J132").Value = 5889.6924
$ws.Range("K132").Value = 4237.7142
$ws.Range("L132").Value = 17669.0772
$ws.Range("M132").Value = -1707.7142
$ws.Range("N132").Value = -22729.0772
$ws.Range("H136").Value = 2746.3171
$ws.Range("I136").Value = 1272.7727
$ws.Range("J136").Value = 4452.5264
$ws.Range("K136").Value = 3818.3181
$ws.Range("L136").Value = 13357.5792
$ws.Range("M136").Value = -1268.3181
$ws.Range("N136").Value = -18457.5792

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H51").Value = 29000
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 29000
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 29000
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -30020
$ws.Range("H122").Value = 72025.64
$ws.Range("I122").Value = 100334.4
$ws.Range("J122").Value = 1253.75
$ws.Range("K122").Value = 301003.2
$ws.Range("L122").Value = 3761.25
$ws.Range("M122").Value = -298553.2
$ws.Range("N122").Value = -8661.25
$ws.Range("H132").Value = 13159949
$ws.Range("I132").Value = 18520200
$ws.Range("J132").Value = 2969.0908
$ws.Range("K132").Value = 55560600
$ws.Range("L132").Value = 8907.2724
$ws.Range("M132").Value = -55558070
$ws.Range("N132").Value = -13967.2724
$ws.Range("H136").Value = 8360345
$ws.Range("I136").Value = 25718634
$ws.Range("J136").Value = 2650.111
$ws.Range("K136").Value = 77155902
$ws.Range("L136").Value = 7950.333
$ws.Range("M136").Value = -77153352
$ws.Range("N136").Value = -13050.333
